# Fixed Bento 80 Test scripts
# Append "order By ... LIMIT 100" clauses to the Neo4j queries stored in
# column B (CasesTab / SamplesTab / FilesTab rows) on the "startup" sheet,
# and move the active selection to C4 (mirrors the saved workbook's
# last-used selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - CasesTab query
$casesQuery = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100 "

# Row 3 - SamplesTab query
$samplesQuery = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"

# Row 4 - FilesTab query
$filesQuery = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = $filesQuery + "`n order By f.file_name ASC LIMIT 100"

# Move selection to C4 (also clears the old frozen/scrolled topLeftCell="A3")
$ws.Range("C4").Select()
